{"js": "// Update the 20x5 arithmetic table (100 cells) with new problem text,\n// in row-major order, matching the order the cells appear in the document.\nconst newValues = [\n  \"20+38=\",\n  \"61-48=\",\n  \"91-81=\",\n  \"37+14=\",\n  \"73-0=\",\n  \"7+60=\",\n  \"8+51=\",\n  \"35+42=\",\n  \"30-19=\",\n  \"59+36=\",\n  \"21-18=\",\n  \"98-42=\",\n  \"44+24=\",\n  \"99-66=\",\n  \"15+20=\",\n  \"64-52=\",\n  \"31-8=\",\n  \"98-8=\",\n  \"69-68=\",\n  \"73-64=\",\n  \"56-8=\",\n  \"82-4=\",\n  \"69-57=\",\n  \"2-0=\",\n  \"82-41=\",\n  \"2+56=\",\n  \"2+69=\",\n  \"53+42=\",\n  \"31+9=\",\n  \"67-41=\",\n  \"65-0=\",\n  \"18+51=\",\n  \"50+24=\",\n  \"72-30=\",\n  \"76-36=\",\n  \"45-36=\",\n  \"58-21=\",\n  \"11-3=\",\n  \"78-30=\",\n  \"23+17=\",\n  \"81-53=\",\n  \"18+36=\",\n  \"71-27=\",\n  \"18+38=\",\n  \"56-39=\",\n  \"65-19=\",\n  \"45+15=\",\n  \"32-1=\",\n  \"62-9=\",\n  \"30+54=\",\n  \"6+15=\",\n  \"74+2=\",\n  \"64-41=\",\n  \"68-49=\",\n  \"2+64=\",\n  \"42-3=\",\n  \"31-1=\",\n  \"93-55=\",\n  \"8+61=\",\n  \"66-38=\",\n  \"30+5=\",\n  \"30+58=\",\n  \"61+10=\",\n  \"49+11=\",\n  \"15-4=\",\n  \"11+15=\",\n  \"76-2=\",\n  \"46+47=\",\n  \"7+82=\",\n  \"56+6=\",\n  \"90-20=\",\n  \"36-20=\",\n  \"39+50=\",\n  \"13+23=\",\n  \"45+22=\",\n  \"51+29=\",\n  \"62-20=\",\n  \"76-33=\",\n  \"24+54=\",\n  \"48+40=\",\n  \"7+90=\",\n  \"45-41=\",\n  \"73-53=\",\n  \"93-31=\",\n  \"95-76=\",\n  \"63-41=\",\n  \"31-18=\",\n  \"1+57=\",\n  \"92-44=\",\n  \"13-5=\",\n  \"86-41=\",\n  \"50+45=\",\n  \"1+57=\",\n  \"40+5=\",\n  \"87+5=\",\n  \"86-59=\",\n  \"8+87=\",\n  \"99-14=\",\n  \"41-31=\",\n  \"16+46=\"\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document, found none.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = 5; // fixed 5-column grid per the document's tblGrid\nlet idx = 0;\n\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    if (idx >= newValues.length) break;\n    const cell = table.getCell(r, c);\n    cell.value = newValues[idx];\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the 20x5 arithmetic table (100 cells) with new problem text,\n# in row-major order, matching the order the cells appear in the document.\n$newValues = @(\n    \"20+38=\",\n    \"61-48=\",\n    \"91-81=\",\n    \"37+14=\",\n    \"73-0=\",\n    \"7+60=\",\n    \"8+51=\",\n    \"35+42=\",\n    \"30-19=\",\n    \"59+36=\",\n    \"21-18=\",\n    \"98-42=\",\n    \"44+24=\",\n    \"99-66=\",\n    \"15+20=\",\n    \"64-52=\",\n    \"31-8=\",\n    \"98-8=\",\n    \"69-68=\",\n    \"73-64=\",\n    \"56-8=\",\n    \"82-4=\",\n    \"69-57=\",\n    \"2-0=\",\n    \"82-41=\",\n    \"2+56=\",\n    \"2+69=\",\n    \"53+42=\",\n    \"31+9=\",\n    \"67-41=\",\n    \"65-0=\",\n    \"18+51=\",\n    \"50+24=\",\n    \"72-30=\",\n    \"76-36=\",\n    \"45-36=\",\n    \"58-21=\",\n    \"11-3=\",\n    \"78-30=\",\n    \"23+17=\",\n    \"81-53=\",\n    \"18+36=\",\n    \"71-27=\",\n    \"18+38=\",\n    \"56-39=\",\n    \"65-19=\",\n    \"45+15=\",\n    \"32-1=\",\n    \"62-9=\",\n    \"30+54=\",\n    \"6+15=\",\n    \"74+2=\",\n    \"64-41=\",\n    \"68-49=\",\n    \"2+64=\",\n    \"42-3=\",\n    \"31-1=\",\n    \"93-55=\",\n    \"8+61=\",\n    \"66-38=\",\n    \"30+5=\",\n    \"30+58=\",\n    \"61+10=\",\n    \"49+11=\",\n    \"15-4=\",\n    \"11+15=\",\n    \"76-2=\",\n    \"46+47=\",\n    \"7+82=\",\n    \"56+6=\",\n    \"90-20=\",\n    \"36-20=\",\n    \"39+50=\",\n    \"13+23=\",\n    \"45+22=\",\n    \"51+29=\",\n    \"62-20=\",\n    \"76-33=\",\n    \"24+54=\",\n    \"48+40=\",\n    \"7+90=\",\n    \"45-41=\",\n    \"73-53=\",\n    \"93-31=\",\n    \"95-76=\",\n    \"63-41=\",\n    \"31-18=\",\n    \"1+57=\",\n    \"92-44=\",\n    \"13-5=\",\n    \"86-41=\",\n    \"50+45=\",\n    \"1+57=\",\n    \"40+5=\",\n    \"87+5=\",\n    \"86-59=\",\n    \"8+87=\",\n    \"99-14=\",\n    \"41-31=\",\n    \"16+46=\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = 5\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    if ($idx -ge $newValues.Count) { break }\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $newValues[$idx]\n    $idx++\n  }\n}\n"}
